$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 13.377
$ws.Range("A8").Value = -21.107
$ws.Range("A10").Value = -20.973
$ws.Range("A12").Value = -21.694
$ws.Range("B13").Value = 6.606
$ws.Range("A18").Value = -21.694
$ws.Range("E20").Value = 12.932
